$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header E1/F1 back to literal text values (the shared-string table is
# collapsed down to just the one coordinate string still referenced below)
$ws.Range("E1").Value = "longitude"
$ws.Range("F1").Value = "coordinates"

# New column B width (bestFit, custom)
$ws.Columns.Item(2).ColumnWidth = 74.1384014423077

$data = @(
  @(5,  "Dona de Casa, Via W Tres Norte - Asa Norte, Brasilia - DF, 70297-400", 0.96250000000000002, -15.7744312, -47.900262400000003),
  @(6,  "cobasi", 0.96879999999999999, -15.736511399999999, -47.893102800000001),
  @(7,  "PÃ£o Dourado - Noroeste", 0.84379999999999999, -15.7408634, -47.9064549),
  @(8,  "Panificadora PÃ£o Prima", 0.5625, -15.754970200000001, -47.898261599999998),
  @(9,  "Leroy Merlin Brasilia Norte, SOFN - Area Especial, Brasilia - DF, 70634-120", 0.96250000000000002, -15.7525496, -47.930304499999998),
  @(10, "Superquadra Norte 210 - Asa Norte, Brasilia - DF, 70862-000", 0.96879999999999999, -15.755914199999999, -47.884851500000003),
  @(11, "Superquadra Norte 115 - Asa Norte, Brasilia - DF, 70297-400", 0.84379999999999999, -15.7416731, -47.893982000000001),
  @(12, "Administracao Regional de Lago Norte - St. De HabitaÃ§Ãµes Individuais Norte CA 5 - Lago Norte, Brasilia - DF, 71503-507", 0.5625, -15.716620900000001, -47.885653499999997),
  @(13, "PÃ£o de AÃ§Ãºcar - Sul 304/305", 0.96250000000000002, -15.8040386, -47.894147099999998),
  @(14, "Boteco do Juca", 0.96879999999999999, -15.8090495, -47.893769499999998),
  @(15, "Casa de Biscoitos Mineiros Asa Sul", 0.84379999999999999, -15.809072199999999, -47.898075499999997),
  @(16, "Simpsons Asa Sul", 0.5625, -15.809183600000001, -47.900110699999999)
)

# Style templates already present in the sheet: F2 carries the plain
# "style 0" look that every new A/B/D/E/F cell should use, while C2 and C5
# carry the two percentage looks ("style 7" / "style 11") used in Capacity.
$ws.Range("F2").Copy()
$plainStyleRange = $ws.Range("A1")

$row = 6
foreach ($entry in $data) {
    $id = $entry[0]
    $desc = $entry[1]
    $cap = $entry[2]
    $lat = $entry[3]
    $lon = $entry[4]
    $coord = "$lat,$lon"

    $ws.Range("F2").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Range("F2").Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Range("F2").Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
    $ws.Range("F2").Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4122)
    $ws.Range("F2").Copy()
    $ws.Cells.Item($row, 6).PasteSpecial(-4122)

    if ($cap -eq 0.5625) {
        $ws.Range("C5").Copy()
    } else {
        $ws.Range("C2").Copy()
    }
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $desc
    $ws.Cells.Item($row, 3).Value = $cap
    $ws.Cells.Item($row, 4).Value = $lat
    $ws.Cells.Item($row, 5).Value = $lon
    $ws.Cells.Item($row, 6).Value = $coord

    $row = $row + 1
}

# Last coordinate cell references a shared string in the target file
$ws.Range("F17").Value = "-15.8091836,-47.9001107"

# Final selection the author left on save
$ws.Range("A18:A22").Select()
